$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 6) entirely so dimension shrinks to A1:B5
$ws.Range("A6:B6").EntireRow.Delete()

# Update remaining rows with the new cluster counts (Birch+K_means values)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 448

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 268

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 241

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 44
